$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H27").Value = 55000
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 55000
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 165000
$ws.Range("M27").Value = $null
$ws.Range("N27").Value = -165202
$ws.Range("H33").Value = 1011957.44
$ws.Range("I33").Value = 1852323
$ws.Range("K33").Value = 1852323
$ws.Range("M33").Value = -1852094
$ws.Range("H40").Value = 5501
$ws.Range("I40").Value = 10000
$ws.Range("K40").Value = 10000
$ws.Range("M40").Value = -9825
$ws.Range("H62").Value = 7859
$ws.Range("I62").Value = 7909.4287
$ws.Range("K62").Value = 7909.4287
$ws.Range("M62").Value = -7285.4287
$ws.Range("H65").Value = 7859
$ws.Range("I65").Value = 7909.4287
$ws.Range("K65").Value = 39547.14350000001
$ws.Range("M65").Value = -36427.14350000001
$ws.Range("H70").Value = 3446.1904
$ws.Range("J70").Value = 3558.5
$ws.Range("L70").Value = 10675.5
$ws.Range("N70").Value = -11215.5
$ws.Range("H73").Value = 3446.1904
$ws.Range("J73").Value = 3558.5
$ws.Range("L73").Value = 10675.5
$ws.Range("N73").Value = -12547.5
$ws.Range("H86").Value = 5042.5713
$ws.Range("I86").Value = 4866.3335
$ws.Range("J86").Value = 6100
$ws.Range("K86").Value = 4866.3335
$ws.Range("L86").Value = 6100
$ws.Range("M86").Value = -3743.3335
$ws.Range("N86").Value = -8346
$ws.Range("H89").Value = 5042.5713
$ws.Range("I89").Value = 4866.3335
$ws.Range("J89").Value = 6100
$ws.Range("K89").Value = 24331.6675
$ws.Range("L89").Value = 30500
$ws.Range("M89").Value = -18715.6675
$ws.Range("N89").Value = -41732
$ws.Range("H98").Value = 1519.7778
$ws.Range("J98").Value = 1190
$ws.Range("L98").Value = 1190
$ws.Range("N98").Value = -4186
$ws.Range("H103").Value = 854.0714
$ws.Range("I103").Value = 311.94446
$ws.Range("J103").Value = 1829.9
$ws.Range("K103").Value = 935.83338
$ws.Range("L103").Value = 5489.700000000001
$ws.Range("M103").Value = -349.83338
$ws.Range("N103").Value = -6661.700000000001
$ws.Range("H106").Value = 6184.1665
$ws.Range("I106").Value = 6633.1113
$ws.Range("J106").Value = 4837.3335
$ws.Range("K106").Value = 6633.1113
$ws.Range("L106").Value = 4837.3335
$ws.Range("M106").Value = -6002.1113
$ws.Range("N106").Value = -6099.3335
$ws.Range("H122").Value = 1519.7778
$ws.Range("J122").Value = 1190
$ws.Range("L122").Value = 3570
$ws.Range("N122").Value = -8470
$ws.Range("H137").Value = 1728.697
$ws.Range("I137").Value = 1454.7727
$ws.Range("J137").Value = 2276.5454
$ws.Range("K137").Value = 4364.3181
$ws.Range("L137").Value = 6829.6362
$ws.Range("M137").Value = -1814.3181
$ws.Range("N137").Value = -11929.6362
$ws.Range("H138").Value = 2924.718
$ws.Range("I138").Value = 2310.0908
$ws.Range("J138").Value = 3166.1785
$ws.Range("K138").Value = 6930.2724
$ws.Range("L138").Value = 9498.5355
$ws.Range("M138").Value = -1790.2724
$ws.Range("N138").Value = -19778.5355

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2176.5356
$ws.Range("I45").Value = 1791.2632
$ws.Range("K45").Value = 1791.2632
$ws.Range("M45").Value = -1414.2632
$ws.Range("H61").Value = 3864.0833
$ws.Range("I61").Value = 3204.6
$ws.Range("K61").Value = 3204.6
$ws.Range("M61").Value = -2992.6
$ws.Range("H132").Value = 3034.5217
$ws.Range("I132").Value = 2146.0789
$ws.Range("K132").Value = 6438.236699999999
$ws.Range("M132").Value = -3908.236699999999
$ws.Range("H136").Value = 3864.0833
$ws.Range("I136").Value = 3204.6
$ws.Range("K136").Value = 9613.799999999999
$ws.Range("M136").Value = -7063.799999999999

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3898.54
$ws.Range("I134").Value = 3565.6052
$ws.Range("K134").Value = 10696.8156
$ws.Range("M134").Value = -8161.8156
$ws.Range("H140").Value = 298944.5
$ws.Range("J140").Value = 371926
$ws.Range("L140").Value = 371926
$ws.Range("N140").Value = -382286

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3281.5715
$ws.Range("I31").Value = 1358.875
$ws.Range("J31").Value = 5845.1665
$ws.Range("K31").Value = 1358.875
$ws.Range("L31").Value = 5845.1665
$ws.Range("M31").Value = -1063.875
$ws.Range("N31").Value = -6435.1665
$ws.Range("H34").Value = 3281.5715
$ws.Range("I34").Value = 1358.875
$ws.Range("J34").Value = 5845.1665
$ws.Range("K34").Value = 1358.875
$ws.Range("L34").Value = 5845.1665
$ws.Range("M34").Value = -1156.875
$ws.Range("N34").Value = -6249.1665
$ws.Range("H62").Value = 3497.1428
$ws.Range("I62").Value = 3462.3333
$ws.Range("J62").Value = 3523.25
$ws.Range("K62").Value = 3462.3333
$ws.Range("L62").Value = 3523.25
$ws.Range("M62").Value = -2838.3333
$ws.Range("N62").Value = -4771.25
$ws.Range("H65").Value = 3497.1428
$ws.Range("I65").Value = 3462.3333
$ws.Range("J65").Value = 3523.25
$ws.Range("K65").Value = 17311.6665
$ws.Range("L65").Value = 17616.25
$ws.Range("M65").Value = -14191.6665
$ws.Range("N65").Value = -23856.25
$ws.Range("H99").Value = 3357.5789
$ws.Range("I99").Value = 3841.6155
$ws.Range("J99").Value = 2308.8333
$ws.Range("K99").Value = 3841.6155
$ws.Range("L99").Value = 2308.8333
$ws.Range("M99").Value = -2343.6155
$ws.Range("N99").Value = -5304.8333
$ws.Range("H105").Value = 2191.0908
$ws.Range("I105").Value = 2226.5
$ws.Range("K105").Value = 2226.5
$ws.Range("M105").Value = -479.5
$ws.Range("H126").Value = 3357.5789
$ws.Range("I126").Value = 3841.6155
$ws.Range("J126").Value = 2308.8333
$ws.Range("K126").Value = 11524.8465
$ws.Range("L126").Value = 6926.499899999999
$ws.Range("M126").Value = -9054.8465
$ws.Range("N126").Value = -11866.4999

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 4429.9473
$ws.Range("I14").Value = 4429.9473
$ws.Range("K14").Value = 13289.8419
$ws.Range("M14").Value = -13116.8419
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").Value = $null
$ws.Range("H121").Value = 22694700
$ws.Range("I121").Value = 10598.5
$ws.Range("J121").Value = 40841980
$ws.Range("K121").Value = 31795.5
$ws.Range("L121").Value = 122525940
$ws.Range("M121").Value = -30485.5
$ws.Range("N121").Value = -122528560
$ws.Range("H131").Value = 1887.4445
$ws.Range("I131").Value = 1174.7273
$ws.Range("J131").Value = 3007.4285
$ws.Range("K131").Value = 3524.1819
$ws.Range("L131").Value = 9022.2855
$ws.Range("M131").Value = 1515.8181
$ws.Range("N131").Value = -19102.2855

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4260.4
$ws.Range("I80").Value = 4195
$ws.Range("K80").Value = 4195
$ws.Range("M80").Value = -3197
$ws.Range("H83").Value = 4260.4
$ws.Range("I83").Value = 4195
$ws.Range("K83").Value = 20975
$ws.Range("M83").Value = -15983
$ws.Range("H95").Value = 116499.836
$ws.Range("J95").Value = 116499.836
$ws.Range("L95").Value = 116499.836
$ws.Range("N95").Value = -121991.836
$ws.Range("H102").Value = 4368.278
$ws.Range("I102").Value = 3972.6667
$ws.Range("K102").Value = 3972.6667
$ws.Range("M102").Value = -2350.6667
$ws.Range("H122").Value = 1784834.4
$ws.Range("I122").Value = 2375479.5
$ws.Range("K122").Value = 7126438.5
$ws.Range("M122").Value = -7123988.5
$ws.Range("H132").Value = 3842.9429
$ws.Range("I132").Value = 4141.7144
$ws.Range("J132").Value = 3394.7856
$ws.Range("K132").Value = 12425.1432
$ws.Range("L132").Value = 10184.3568
$ws.Range("M132").Value = -9895.143199999999
$ws.Range("N132").Value = -15244.3568

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1865.0769
$ws.Range("I82").Value = 1920.5
$ws.Range("K82").Value = 1920.5
$ws.Range("M82").Value = -1559.5
$ws.Range("H85").Value = 1865.0769
$ws.Range("I85").Value = 1920.5
$ws.Range("K85").Value = 1920.5
$ws.Range("M85").Value = -672.5
$ws.Range("H93").Value = 3440
$ws.Range("I93").Value = 3575.125
$ws.Range("K93").Value = 3575.125
$ws.Range("M93").Value = -2327.125
$ws.Range("H136").Value = 10757049
$ws.Range("I136").Value = 2672.6316
$ws.Range("J136").Value = 27784812
$ws.Range("K136").Value = 8017.8948
$ws.Range("L136").Value = 83354436
$ws.Range("M136").Value = -5467.8948
$ws.Range("N136").Value = -83359536

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 40619.25
$ws.Range("J105").Value = 40619.25
$ws.Range("L105").Value = 40619.25
$ws.Range("N105").Value = -47607.25
$ws.Range("H132").Value = 5160.909
$ws.Range("I132").Value = 2781.08
$ws.Range("J132").Value = 12597.875
$ws.Range("K132").Value = 8343.24
$ws.Range("L132").Value = 37793.625
$ws.Range("M132").Value = -5813.24
$ws.Range("N132").Value = -42853.625
